$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '57.229.55'
$ws.Range("E2").Value2 = '  -1.65%  '
$ws.Range("D3").Value2 = '3.081.03'
$ws.Range("E3").Value2 = '  -1.34%  '
$ws.Range("E4").Value2 = '  -0.06%  '
$ws.Range("D5").Value2 = '''522.88'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value2 = '  -1.01%  '
$ws.Range("D6").Value2 = '''135.72'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value2 = '  -4.78%  '
$ws.Range("E7").Value2 = '  -0.08%  '
$ws.Range("D8").Value2 = '3.077.85'
$ws.Range("E8").Value2 = '  -1.39%  '
$ws.Range("D9").Value2 = '''0.464'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value2 = '  +4.36%  '
$ws.Range("D10").Value2 = '''7.34'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value2 = '  +2.25%  '
$ws.Range("E11").Value2 = '  -2.25%  '
$ws.Range("D12").Value2 = '''0.400'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value2 = '  +1.44%  '
$ws.Range("D13").Value2 = '''0.136'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value2 = '  +1.24%  '
$ws.Range("D14").Value2 = '3.601.57'
$ws.Range("E14").Value2 = '  -1.53%  '
$ws.Range("D15").Value2 = '''25.24'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value2 = '  -2.32%  '
$ws.Range("E16").Value2 = '  -2.66%  '
$ws.Range("D17").Value2 = '57.269.89'
$ws.Range("E17").Value2 = '  -1.67%  '
$ws.Range("D18").Value2 = '3.074.38'
$ws.Range("E18").Value2 = '  -1.43%  '
$ws.Range("D19").Value2 = '''5.87'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value2 = '  -4.11%  '
$ws.Range("D20").Value2 = '''12.46'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value2 = '  -2.74%  '
$ws.Range("D21").Value2 = '''7.83'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value2 = '  -2.11%  '
$ws.Range("D22").Value2 = '''350.80'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value2 = '  +2.08%  '
$ws.Range("E23").Value2 = '  +0.35%  '
$ws.Range("D24").Value2 = '''69.02'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value2 = '  +1.78%  '
$ws.Range("D25").Value2 = '''0.498'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value2 = '  -3.21%  '
$ws.Range("D27").Value2 = '''0.998'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value2 = '  -0.19%  '
$ws.Range("D28").Value2 = '0.0₃0868'
$ws.Range("E28").Value2 = '  -6.42%  '
$ws.Range("D29").Value2 = '''1.00'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value2 = '  +0.08%  '
$ws.Range("D30").Value2 = '''7.21'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value2 = '  -1.41%  '
$ws.Range("E31").Value2 = '  -1.09%  '
$ws.Range("D32").Value2 = '''5.85'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value2 = '  -8.35%  '
$ws.Range("D33").Value2 = '''20.96'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value2 = '  -0.57%  '
$ws.Range("B34").Value2 = 'NEARProtocol'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value2 = '''4.82'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value2 = '  +3.12%  '
$ws.Range("B35").Value2 = 'Monero'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value2 = '''158.83'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value2 = '  +0.45%  '
$ws.Range("D36").Value2 = '''1.13'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value2 = '  -5.23%  '
$ws.Range("D37").Value2 = '''6.00'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value2 = '  -3.43%  '
$ws.Range("D38").Value2 = '''25.64'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value2 = '  -2.70%  '
$ws.Range("D39").Value2 = '''1.23'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value2 = '  -2.26%  '
$ws.Range("E40").Value2 = '  -2.18%  '
$ws.Range("E41").Value2 = '  -3.76%  '
$ws.Range("D42").Value2 = '''4.07'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value2 = '  +1.45%  '
$ws.Range("D43").Value2 = '''0.692'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value2 = '  -0.09%  '
$ws.Range("D44").Value2 = '2.400.09'
$ws.Range("E44").Value2 = '  +5.33%  '
$ws.Range("D45").Value2 = '''36.67'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value2 = '  +0.26%  '
$ws.Range("E46").Value2 = '  -0.02%  '
$ws.Range("D47").Value2 = '3.117.88'
$ws.Range("E47").Value2 = '  -1.34%  '
$ws.Range("D48").Value2 = '''0.0261'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value2 = '  -0.84%  '
$ws.Range("D49").Value2 = '''0.948'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value2 = '  -5.63%  '
$ws.Range("D50").Value2 = '''5.97'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value2 = '  -2.48%  '
$ws.Range("D51").Value2 = '''19.63'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value2 = '  -5.08%  '
